$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the existing hyperlink addresses keyed by their (pre-shift) row,
# since the underlying engine does not relocate Hyperlink.Range when the
# sheet is shifted by a column insert.
$hyperlinkAddresses = @{}
foreach ($hl in $ws.Hyperlinks) {
    $hyperlinkAddresses[$hl.Range.Row] = $hl.Address
}

# Shift everything from columns A:B to B:C by inserting a blank column
# before column A (equivalent of cutting A:B and pasting into B:C).
$ws.Columns("A:A").Insert()

# Drop the stale hyperlinks (still anchored to column A) and recreate them
# one column to the right, on top of the data that has already shifted.
$ws.Hyperlinks.Delete()
foreach ($row in $hyperlinkAddresses.Keys) {
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 2), $hyperlinkAddresses[$row]) | Out-Null
}

# Match the saved selection state.
$ws.Range("B11").Select()
